$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb2"
$ws.Range("C2").Value = "Tgfbr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.465615333333333
$ws.Range("H2").Value = 4.396846
$ws.Range("I2").Value = 0.04672291954663727
$ws.Range("J2").Value = 0.04672291954663728
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 68.65869266666667
$ws.Range("N2").Value = 205.976078
$ws.Range("O2").Value = 0.6475952735309433
$ws.Range("P2").Value = 0.6475952735309431
$ws.Range("Q2").Value = 100.6272327388876
$ws.Range("R2").Value = 905.645094649988
$ws.Range("S2").Value = 0.03025754186396882
$ws.Range("T2").Value = 0.03025754186396882

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb2"
$ws.Range("C3").Value = "Tgfbr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.465615333333333
$ws.Range("H3").Value = 4.396846
$ws.Range("I3").Value = 0.04672291954663727
$ws.Range("J3").Value = 0.04672291954663728
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.23943666666667
$ws.Range("N3").Value = 84.71831
$ws.Range("O3").Value = 0.2663570336431459
$ws.Range("P3").Value = 0.2663570336431459
$ws.Range("Q3").Value = 41.38815138336222
$ws.Range("R3").Value = 372.49336245026
$ws.Range("S3").Value = 0.01244497825358966
$ws.Range("T3").Value = 0.01244497825358966

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb2"
$ws.Range("C4").Value = "Tgfbr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.465615333333333
$ws.Range("H4").Value = 4.396846
$ws.Range("I4").Value = 0.04672291954663727
$ws.Range("J4").Value = 0.04672291954663728
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.122861666666667
$ws.Range("N4").Value = 27.368585
$ws.Range("O4").Value = 0.08604769282591093
$ws.Range("P4").Value = 0.08604769282591092
$ws.Range("Q4").Value = 13.37060594254556
$ws.Range("R4").Value = 120.33545348291
$ws.Range("S4").Value = 0.004020399429078793
$ws.Range("T4").Value = 0.004020399429078793

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tgfb2"
$ws.Range("C5").Value = "Tgfbr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.88237266666667
$ws.Range("H5").Value = 56.64711800000001
$ws.Range("I5").Value = 0.6019584804341267
$ws.Range("J5").Value = 0.6019584804341268
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 68.65869266666667
$ws.Range("N5").Value = 205.976078
$ws.Range("O5").Value = 0.6475952735309433
$ws.Range("P5").Value = 0.6475952735309431
$ws.Range("Q5").Value = 1296.439021738134
$ws.Range("R5").Value = 11667.95119564321
$ws.Range("S5").Value = 0.3898254667910092
$ws.Range("T5").Value = 0.3898254667910093

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb2"
$ws.Range("C6").Value = "Tgfbr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 18.88237266666667
$ws.Range("H6").Value = 56.64711800000001
$ws.Range("I6").Value = 0.6019584804341267
$ws.Range("J6").Value = 0.6019584804341268
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 28.23943666666667
$ws.Range("N6").Value = 84.71831
$ws.Range("O6").Value = 0.2663570336431459
$ws.Range("P6").Value = 0.2663570336431459
$ws.Range("Q6").Value = 533.2275670367312
$ws.Range("R6").Value = 4799.04810333058
$ws.Range("S6").Value = 0.1603358752247697
$ws.Range("T6").Value = 0.1603358752247697

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb2"
$ws.Range("C7").Value = "Tgfbr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 18.88237266666667
$ws.Range("H7").Value = 56.64711800000001
$ws.Range("I7").Value = 0.6019584804341267
$ws.Range("J7").Value = 0.6019584804341268
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.122861666666667
$ws.Range("N7").Value = 27.368585
$ws.Range("O7").Value = 0.08604769282591093
$ws.Range("P7").Value = 0.08604769282591092
$ws.Range("Q7").Value = 172.2612737764478
$ws.Range("R7").Value = 1550.35146398803
$ws.Range("S7").Value = 0.05179713841834785
$ws.Range("T7").Value = 0.05179713841834785

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tgfb2"
$ws.Range("C8").Value = "Tgfbr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.020243
$ws.Range("H8").Value = 33.060729
$ws.Range("I8").Value = 0.351318600019236
$ws.Range("J8").Value = 0.351318600019236
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 68.65869266666667
$ws.Range("N8").Value = 205.976078
$ws.Range("O8").Value = 0.6475952735309433
$ws.Range("P8").Value = 0.6475952735309431
$ws.Range("Q8").Value = 756.6354772489847
$ws.Range("R8").Value = 6809.719295240862
$ws.Range("S8").Value = 0.2275122648759652
$ws.Range("T8").Value = 0.2275122648759652

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tgfb2"
$ws.Range("C9").Value = "Tgfbr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.020243
$ws.Range("H9").Value = 33.060729
$ws.Range("I9").Value = 0.351318600019236
$ws.Range("J9").Value = 0.351318600019236
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.23943666666667
$ws.Range("N9").Value = 84.71831
$ws.Range("O9").Value = 0.2663570336431459
$ws.Range("P9").Value = 0.2663570336431459
$ws.Range("Q9").Value = 311.2054542497767
$ws.Range("R9").Value = 2800.84908824799
$ws.Range("S9").Value = 0.09357618016478657
$ws.Range("T9").Value = 0.09357618016478657

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb2"
$ws.Range("C10").Value = "Tgfbr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.020243
$ws.Range("H10").Value = 33.060729
$ws.Range("I10").Value = 0.351318600019236
$ws.Range("J10").Value = 0.351318600019236
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.122861666666667
$ws.Range("N10").Value = 27.368585
$ws.Range("O10").Value = 0.08604769282591093
$ws.Range("P10").Value = 0.08604769282591092
$ws.Range("Q10").Value = 100.5361524220517
$ws.Range("R10").Value = 904.8253717984651
$ws.Range("S10").Value = 0.03023015497848429
$ws.Range("T10").Value = 0.03023015497848429
